$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175718784332275
$ws.Range("B1").Value = 2.407128810882568
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.344990015029907
$ws.Range("E1").Value = 1.206812858581543
